$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks before re-arranging columns; they will be
# re-created (pointing at the new "client_email" column, D) further below.
$ws.Hyperlinks.Delete()

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "client_name"
$ws.Range("B1").Value = "company_name"
$ws.Range("C1").Value = "phone_number"
$ws.Range("D1").Value = "client_email"
$ws.Range("E1").Value = "customer_type"

# --- Existing client rows, now spread across the new column layout -------
$ws.Range("A2").Value = "Rachit"
$ws.Range("D2").Value = "rachitarora1993@gmail.com"

$ws.Range("A3").Value = "IIT Rachit"
$ws.Range("D3").Value = "rachitar@iitrpr.ac.in"

$ws.Range("A4").Value = "Aman"
$ws.Range("D4").Value = "sikarwaraman26@gmail.com"

# --- New row (Rishabh submission) -----------------------------------------
$ws.Range("A5").Value = "Rishabh"
$ws.Range("D5").Value = "rishwebd@gmail.com"

# --- Fill in the new company_name / phone_number / customer_type columns -
# (B2:B4 previously held the hyperlinked e-mail column, so reset their style
# back to Normal before reusing them for plain text.)
$ws.Range("B2:B4").Style = "Normal"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = "OM"
    $ws.Cells.Item($r, 3).Value = "1234-567-890"
    $ws.Cells.Item($r, 5).Value = "B2B"
}

# --- Re-create the mailto hyperlinks on column D --------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:rachitarora1993@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:rachitar@iitrpr.ac.in")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:sikarwaraman26@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:rishwebd@gmail.com")
$ws.Range("D2:D5").Style = "Hyperlink"

# --- Column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 31.5
$ws.Columns.Item(2).ColumnWidth = 46.666666666666664
$ws.Columns.Item(4).ColumnWidth = 34.5
$ws.Columns.Item(5).ColumnWidth = 36.666666666666664

# --- Selection / active cell ----------------------------------------------
$ws.Range("E2:E5").Select() | Out-Null
